# se modif data para Smoke en QA
# Update the "Smoke" test data (account, motor) from the previous run (SMA017 / siete)
# to the new QA run (QADos / SMA018), bump the running counters, and move the
# remembered cell selection on each sheet down one row (as Excel does after typing
# into the last used row and pressing Enter).

$wb = $excel.ActiveWorkbook

# --- DatosCuenta -----------------------------------------------------------
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokQADos"
$wsCuenta.Range("B2").Value = "SmokeNameQADos"
$wsCuenta.Range("C2").Value = 27100115
$wsCuenta.Range("D2").Value = 117
$wsCuenta.Range("D3").Select()

# --- DatosHogar --------------------------------------------------------------
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 636
$wsHogar.Range("A3").Select()

# --- DatosMotor --------------------------------------------------------------
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMA018"
$wsMotor.Range("B2").Value = "ABC12SSMA018"
$wsMotor.Range("C2").Value = "ZAZ123SSMA018"
$wsMotor.Range("A2:C2").Select()

# --- DatosAP -------------------------------------------------------------------
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200118
$wsAP.Range("A3").Select()
